$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings)
$ws.Range("B1").Value = "A_C"
$ws.Range("C1").Value = "LF_FFR"

# Update numeric values
$ws.Range("B2").Value = 1.061839086070209
$ws.Range("C2").Value = -0.9000230269908649
$ws.Range("C3").Value = 0
